$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Placas")

$ws.Range("A2").Value = "RNH0D65"
$ws.Range("B2").Value = "Diarista 6 - 12"

$ws.Range("A3").ClearContents()
$ws.Range("B3").Clear()

$ws.Activate() | Out-Null
$ws.Range("B5").Select() | Out-Null
